# Apply the "Ajout des nouveaux profils au flux 1" edits to the
# StructureDefinition-tddui-attachment workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": simple Property / Value table ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-ballot -> 2.0.0
$meta.Range("B3").Value = "2.0.0"

# Title: Pièce jointe de l'évènement -> Pièce jointe
$meta.Range("B5").Value = "Pièce jointe"

# Date: 2025-10-01T08:29:05+00:00 -> 2025-10-20T13:10:23+00:00
$meta.Range("B8").Value = "2025-10-20T13:10:23+00:00"

# Description: updated wording mentioning the new profile reference
$meta.Range("B12").Value = "Pièces jointes liées à l’événement et à l'évaluation. L'extension référence le profil PDSm_SimplifiedPublish."

# --- Sheet "Elements": element table ---
$elements = $wb.Worksheets.Item("Elements")

# Row 6 is "Extension.value[x]"; column K is "Type(s)".
# Attachment -> Reference(.../tddui-document-reference)
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-document-reference)
"

# Widen column K (Type(s)) to fit the new, longer reference text.
# (COM ColumnWidth is in characters; the engine stores/rounds the OOXML
# width to the nearest 1/6, so 78.8333... lands on the closest
# representable value to the target 79.70703125.)
$elements.Columns.Item(11).ColumnWidth = 78.8333333333333
